# Auto-generated Excel COM-interop script
# Applies the 2025-09-11 daily crime data update to violent-crime-full-year.xlsx
# For each affected worksheet, updates the specific cells in column L (2025 YTD totals)
# and a handful of prior-year correction cells, per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4678
$ws.Range("L3").Value = 5051
$ws.Range("H4").Value = 1760
$ws.Range("J4").Value = 1873
$ws.Range("K4").Value = 1781
$ws.Range("L4").Value = 1249
$ws.Range("L6").Value = 4294
$ws.Range("H7").Value = 26076
$ws.Range("J7").Value = 29349
$ws.Range("K7").Value = 27573
$ws.Range("L7").Value = 15568

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 346
$ws.Range("L4").Value = 81
$ws.Range("L6").Value = 273
$ws.Range("L7").Value = 1034

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 137
$ws.Range("L7").Value = 344

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 218
$ws.Range("L7").Value = 708

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 82
$ws.Range("L3").Value = 68
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 215

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 74
$ws.Range("L3").Value = 110
$ws.Range("L7").Value = 269

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 129
$ws.Range("L7").Value = 508
$ws.Range("L8").Value = 1034
$ws.Range("L11").Value = 252
$ws.Range("L18").Value = 109
$ws.Range("L23").Value = 170
$ws.Range("L25").Value = 92
$ws.Range("L29").Value = 858
$ws.Range("L30").Value = 73
$ws.Range("L31").Value = 155
$ws.Range("L33").Value = 708
$ws.Range("L40").Value = 41
$ws.Range("L41").Value = 70
$ws.Range("L42").Value = 509
$ws.Range("L43").Value = 114
$ws.Range("L45").Value = 29
$ws.Range("L48").Value = 203
$ws.Range("L51").Value = 192
$ws.Range("L54").Value = 326
$ws.Range("L55").Value = 148
$ws.Range("H63").Value = 311
$ws.Range("J63").Value = 225
$ws.Range("K63").Value = 168
$ws.Range("L63").Value = 46
$ws.Range("K67").Value = 1071
$ws.Range("L67").Value = 535
$ws.Range("L70").Value = 26
$ws.Range("L76").Value = 244
$ws.Range("L78").Value = 209
$ws.Range("L79").Value = 412
$ws.Range("L83").Value = 344
$ws.Range("L85").Value = 798
$ws.Range("L89").Value = 224
$ws.Range("L91").Value = 209
$ws.Range("L95").Value = 215
$ws.Range("L99").Value = 269
$ws.Range("H101").Value = 26076
$ws.Range("J101").Value = 29349
$ws.Range("K101").Value = 27573
$ws.Range("L101").Value = 15568

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K4").Value = 60
$ws.Range("L6").Value = 123
$ws.Range("K7").Value = 1071
$ws.Range("L7").Value = 535

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 326

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 324
$ws.Range("L7").Value = 858

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 50
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 146
$ws.Range("L3").Value = 171
$ws.Range("L6").Value = 143
$ws.Range("L7").Value = 509

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 68
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 148

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 63
$ws.Range("L7").Value = 170

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 92
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 148
$ws.Range("L7").Value = 412

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 168
$ws.Range("L7").Value = 508

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 93
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 252

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 66
$ws.Range("L4").Value = 36
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 192

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 325
$ws.Range("L7").Value = 798

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 41
